# Apply "model 3 output" values to sheets y1, y2, y3 (A-column only; B column unchanged)
$wb = $excel.ActiveWorkbook

# --- Sheet y1 ---
$ws1 = $wb.Worksheets.Item("y1")
$ws1.Range("A3").Value = 0
$ws1.Range("A4").Value = 1

# --- Sheet y2 ---
$ws2 = $wb.Worksheets.Item("y2")
$ws2.Range("A3").Value = 1
$ws2.Range("A4").Value = 0
$ws2.Range("A5").Value = 1
$ws2.Range("A9").Value = 1

# --- Sheet y3 ---
$ws3 = $wb.Worksheets.Item("y3")
$ws3.Range("A5").Value = 0
$ws3.Range("A9").Value = 0
